# "Generate Report for Handoff"
#
# The localization-status report tracks, per source file and target
# language, the handoff/handback lifecycle. The ab33e234-...md file just
# had a new handoff generated for both target languages (zh-cn, de-de):
#   - Its Status moves from "Handed back: in sync with en-US" to
#     "Ready for handoff" (on the Overview sheet, and on each language
#     sheet's Status column).
#   - Each language sheet's "Latest Handoff Datetime" column is stamped
#     with the new handoff timestamp for that language.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the ab33e234-c250-4f33-9a6d-87ebf8d2dd54.md entry
# (columns: B = zh-cn status, C = de-de status).
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the same entry.
# Column B = Status, Column D = Latest Handoff Datetime.
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-28 04:06:32"

# de-de sheet: row 3 is the same entry.
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-28 04:06:42"
